$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1765.7333
$ws.Range("I11").Value = 1765.7333
$ws.Range("K11").Value = 1765.7333
$ws.Range("M11").Value = -1625.7333

$ws.Range("H74").Value = 93758940
$ws.Range("I74").Value = 214288640
$ws.Range("J74").Value = 13627.667
$ws.Range("K74").Value = 214288640
$ws.Range("L74").Value = 13627.667
$ws.Range("M74").Value = -214287704
$ws.Range("N74").Value = -15499.667

$ws.Range("H77").Value = 93758940
$ws.Range("I77").Value = 214288640
$ws.Range("J77").Value = 13627.667
$ws.Range("K77").Value = 1071443200
$ws.Range("L77").Value = 68138.33499999999
$ws.Range("M77").Value = -1071438520
$ws.Range("N77").Value = -77498.33499999999

$ws.Range("H92").Value = 142858140
$ws.Range("I92").Value = 989.25
$ws.Range("J92").Value = 333334340
$ws.Range("K92").Value = 989.25
$ws.Range("L92").Value = 333334340
$ws.Range("M92").Value = 258.75
$ws.Range("N92").Value = -333336836

$ws.Range("H98").Value = 58828516
$ws.Range("I98").Value = 62504956
$ws.Range("J98").Value = 5500
$ws.Range("K98").Value = 62504956
$ws.Range("L98").Value = 5500
$ws.Range("M98").Value = -62503458
$ws.Range("N98").Value = -8496

$ws.Range("H122").Value = 58828516
$ws.Range("I122").Value = 62504956
$ws.Range("J122").Value = 5500
$ws.Range("K122").Value = 187514868
$ws.Range("L122").Value = 16500
$ws.Range("M122").Value = -187512418
$ws.Range("N122").Value = -21400

$ws.Range("H132").Value = 2701.28
$ws.Range("I132").Value = 2756.6
$ws.Range("K132").Value = 8269.799999999999
$ws.Range("M132").Value = -5739.799999999999

$ws.Range("H135").Value = 1112073.9
$ws.Range("I135").Value = 1429232.8
$ws.Range("K135").Value = 12863095.2
$ws.Range("M135").Value = -12860560.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1393182.6
$ws.Range("I32").Value = 1473405
$ws.Range("J32").Value = 29401.6
$ws.Range("K32").Value = 1473405
$ws.Range("L32").Value = 29401.6
$ws.Range("M32").Value = -1473118
$ws.Range("N32").Value = -29975.6

$ws.Range("H46").Value = 3811.8333
$ws.Range("I46").Value = 2999
$ws.Range("J46").Value = 3974.4
$ws.Range("K46").Value = 2999
$ws.Range("L46").Value = 3974.4
$ws.Range("M46").Value = -2680
$ws.Range("N46").Value = -4612.4

$ws.Range("H61").Value = 7802.0625
$ws.Range("I61").Value = 4259.778
$ws.Range("K61").Value = 4259.778
$ws.Range("M61").Value = -4047.778

$ws.Range("H97").Value = 5209311.5
$ws.Range("I97").Value = 903.5714
$ws.Range("K97").Value = 903.5714
$ws.Range("M97").Value = -407.5714

$ws.Range("H122").Value = 4210.1177
$ws.Range("I122").Value = 2953.182
$ws.Range("K122").Value = 8859.545999999998
$ws.Range("M122").Value = -6409.545999999998

$ws.Range("H132").Value = 4022.8813
$ws.Range("I132").Value = 1869
$ws.Range("K132").Value = 5607
$ws.Range("M132").Value = -3077

$ws.Range("H135").Value = 62695
$ws.Range("J135").Value = 62695
$ws.Range("L135").Value = 62695
$ws.Range("N135").Value = -72835

$ws.Range("H136").Value = 7802.0625
$ws.Range("I136").Value = 4259.778
$ws.Range("K136").Value = 12779.334
$ws.Range("M136").Value = -10229.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 383
$ws.Range("I16").Value = 420
$ws.Range("J16").Value = 309
$ws.Range("K16").Value = 420
$ws.Range("L16").Value = 309
$ws.Range("M16").Value = -250
$ws.Range("N16").Value = -649

$ws.Range("H20").Value = 4763795.5
$ws.Range("I20").Value = 6946382
$ws.Range("K20").Value = 6946382
$ws.Range("M20").Value = -6946135

$ws.Range("H26").Value = 28142
$ws.Range("I26").Value = 22982
$ws.Range("K26").Value = 22982
$ws.Range("M26").Value = -22690

$ws.Range("H96").Value = 22802
$ws.Range("I96").Value = 7249.75
$ws.Range("K96").Value = 7249.75
$ws.Range("M96").Value = -4503.75

$ws.Range("H107").Value = 187506000
$ws.Range("I107").Value = 281252000
$ws.Range("J107").Value = 14000
$ws.Range("K107").Value = 281252000
$ws.Range("L107").Value = 14000
$ws.Range("M107").Value = -281250080
$ws.Range("N107").Value = -17840

$ws.Range("H134").Value = 5552.5293
$ws.Range("I134").Value = 2367.25
$ws.Range("J134").Value = 9430.261
$ws.Range("K134").Value = 7101.75
$ws.Range("L134").Value = 28290.783
$ws.Range("M134").Value = -4566.75
$ws.Range("N134").Value = -33360.783

$ws.Range("H135").Value = 99390
$ws.Range("J135").Value = 99390
$ws.Range("L135").Value = 99390
$ws.Range("N135").Value = -109530

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6014.522
$ws.Range("I16").Value = 4084.3
$ws.Range("J16").Value = 7499.3076
$ws.Range("K16").Value = 4084.3
$ws.Range("L16").Value = 7499.3076
$ws.Range("M16").Value = -3797.3
$ws.Range("N16").Value = -8073.3076

$ws.Range("H31").Value = 7311.404
$ws.Range("I31").Value = 3318.516
$ws.Range("J31").Value = 13205.667
$ws.Range("K31").Value = 3318.516
$ws.Range("L31").Value = 13205.667
$ws.Range("M31").Value = -3023.516
$ws.Range("N31").Value = -13795.667

$ws.Range("H34").Value = 7311.404
$ws.Range("I34").Value = 3318.516
$ws.Range("J34").Value = 13205.667
$ws.Range("K34").Value = 3318.516
$ws.Range("L34").Value = 13205.667
$ws.Range("M34").Value = -3116.516
$ws.Range("N34").Value = -13609.667

$ws.Range("H42").Value = 46166.332
$ws.Range("J42").Value = 44250
$ws.Range("L42").Value = 44250
$ws.Range("N42").Value = -45436

$ws.Range("H58").Value = 9264865
$ws.Range("I58").Value = 16131414
$ws.Range("J58").Value = 9952.434999999999
$ws.Range("K58").Value = 16131414
$ws.Range("L58").Value = 9952.434999999999
$ws.Range("M58").Value = -16131211
$ws.Range("N58").Value = -10358.435

$ws.Range("H99").Value = 4710.1562
$ws.Range("I99").Value = 3053.1333
$ws.Range("K99").Value = 3053.1333
$ws.Range("M99").Value = -1555.1333

$ws.Range("H103").Value = 62627.332
$ws.Range("I103").Value = 53941
$ws.Range("K103").Value = 53941
$ws.Range("M103").Value = -52769

$ws.Range("H105").Value = 4204616.5
$ws.Range("I105").Value = 4763632
$ws.Range("J105").Value = 12000
$ws.Range("K105").Value = 4763632
$ws.Range("L105").Value = 12000
$ws.Range("M105").Value = -4761885
$ws.Range("N105").Value = -15494

$ws.Range("H107").Value = 2782.3157
$ws.Range("I107").Value = 2714.625
$ws.Range("J107").Value = 2831.5454
$ws.Range("K107").Value = 2714.625
$ws.Range("L107").Value = 2831.5454
$ws.Range("M107").Value = -794.625
$ws.Range("N107").Value = -6671.5454

$ws.Range("H113").Value = 6014.522
$ws.Range("I113").Value = 4084.3
$ws.Range("J113").Value = 7499.3076
$ws.Range("K113").Value = 4084.3
$ws.Range("L113").Value = 7499.3076
$ws.Range("M113").Value = -1914.3
$ws.Range("N113").Value = -11839.3076

$ws.Range("H126").Value = 4710.1562
$ws.Range("I126").Value = 3053.1333
$ws.Range("K126").Value = 9159.3999
$ws.Range("M126").Value = -6689.3999

$ws.Range("H132").Value = 4405.629
$ws.Range("I132").Value = 2100.262
$ws.Range("K132").Value = 6300.786
$ws.Range("M132").Value = -3770.786

$ws.Range("H136").Value = 9264865
$ws.Range("I136").Value = 16131414
$ws.Range("J136").Value = 9952.434999999999
$ws.Range("K136").Value = 48394242
$ws.Range("L136").Value = 29857.305
$ws.Range("M136").Value = -48391692
$ws.Range("N136").Value = -34957.305

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 22019632
$ws.Range("I4").Value = 26928770
$ws.Range("J4").Value = 4837654.5
$ws.Range("K4").Value = 80786310
$ws.Range("L4").Value = 14512963.5
$ws.Range("M4").Value = -80786198
$ws.Range("N4").Value = -14513187.5

$ws.Range("H26").Value = 461.8095
$ws.Range("I26").Value = 166
$ws.Range("J26").Value = 580.13336
$ws.Range("K26").Value = 498
$ws.Range("L26").Value = 1740.40008
$ws.Range("M26").Value = -210
$ws.Range("N26").Value = -2316.40008

$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws.Range("H74").Value = 2250
$ws.Range("I74").Value = 2250
$ws.Range("K74").Value = 6750
$ws.Range("M74").Value = -5689

$ws.Range("J75").Value = 16669663
$ws.Range("L75").Value = 50008989
$ws.Range("N75").Value = -50010985

$ws.Range("H76").Value = 3013
$ws.Range("I76").Value = 3013
$ws.Range("K76").Value = 9039
$ws.Range("M76").Value = -8656

$ws.Range("H77").Value = 2250
$ws.Range("I77").Value = 2250
$ws.Range("K77").Value = 20250
$ws.Range("M77").Value = -14946

$ws.Range("J78").Value = 16669663
$ws.Range("L78").Value = 150026967
$ws.Range("N78").Value = -150036951

$ws.Range("H79").Value = 3013
$ws.Range("I79").Value = 3013
$ws.Range("K79").Value = 9039
$ws.Range("M79").Value = -7713

$ws.Range("H92").Value = 7694057
$ws.Range("J92").Value = 7694057
$ws.Range("L92").Value = 23082171
$ws.Range("N92").Value = -23084667

$ws.Range("H131").Value = 2138.04
$ws.Range("J131").Value = 2227.568
$ws.Range("L131").Value = 6682.704000000001
$ws.Range("N131").Value = -16762.704

$ws.Range("H132").Value = 10982.435
$ws.Range("I132").Value = 2518.6875
$ws.Range("J132").Value = 30328.143
$ws.Range("K132").Value = 22668.1875
$ws.Range("L132").Value = 272953.287
$ws.Range("M132").Value = -20138.1875
$ws.Range("N132").Value = -278013.287

$ws.Range("H134").Value = 6828.8
$ws.Range("I134").Value = 5243.2
$ws.Range("K134").Value = 15729.6
$ws.Range("M134").Value = -10659.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 79721.766
$ws.Range("I80").Value = 2360.375
$ws.Range("K80").Value = 2360.375
$ws.Range("M80").Value = -1362.375

$ws.Range("H83").Value = 79721.766
$ws.Range("I83").Value = 2360.375
$ws.Range("K83").Value = 11801.875
$ws.Range("M83").Value = -6809.875

$ws.Range("H102").Value = 1906.2325
$ws.Range("I102").Value = 1588.7894
$ws.Range("K102").Value = 1588.7894
$ws.Range("M102").Value = 33.21060000000011

$ws.Range("H122").Value = 17958642
$ws.Range("I122").Value = 23942854
$ws.Range("K122").Value = 71828562
$ws.Range("M122").Value = -71826112

$ws.Range("H132").Value = 5089.4463
$ws.Range("I132").Value = 3570.3257
$ws.Range("J132").Value = 10114.23
$ws.Range("K132").Value = 10710.9771
$ws.Range("L132").Value = 30342.69
$ws.Range("M132").Value = -8180.9771
$ws.Range("N132").Value = -35402.69

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6947755
$ws.Range("I46").Value = 1432.1428
$ws.Range("K46").Value = 1432.1428
$ws.Range("M46").Value = -1244.1428

$ws.Range("H61").Value = 3938.625
$ws.Range("I61").Value = 1752.6111
$ws.Range("K61").Value = 1752.6111
$ws.Range("M61").Value = -1550.6111

$ws.Range("H100").Value = 3165.543
$ws.Range("I100").Value = 2904.9333
$ws.Range("K100").Value = 2904.9333
$ws.Range("M100").Value = -2363.9333

$ws.Range("H113").Value = 3938.625
$ws.Range("I113").Value = 1752.6111
$ws.Range("K113").Value = 1752.6111
$ws.Range("M113").Value = 417.3888999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws.Range("H132").Value = 23811832
$ws.Range("I132").Value = 26317976
$ws.Range("J132").Value = 3470.625
$ws.Range("K132").Value = 78953928
$ws.Range("L132").Value = 10411.875
$ws.Range("M132").Value = -78951398
$ws.Range("N132").Value = -15471.875

$ws.Range("H136").Value = 23283232
$ws.Range("I136").Value = 40001256
$ws.Range("K136").Value = 120003768
$ws.Range("M136").Value = -120001218
